# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 0b6255bb... / ea522373... handback rows, for both
# the zh-cn and de-de report sheets (rows 21 and 22).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D21").Value = "2016-03-08 02:45:36"
$wsZh.Range("D22").Value = "2016-03-08 02:45:36"
$wsZh.Range("G21").Value = "2016-03-08 02:46:19"
$wsZh.Range("G22").Value = "2016-03-08 02:46:19"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D21").Value = "2016-03-08 02:45:44"
$wsDe.Range("D22").Value = "2016-03-08 02:45:44"
$wsDe.Range("G21").Value = "2016-03-08 02:46:33"
$wsDe.Range("G22").Value = "2016-03-08 02:46:33"
